# "use expriment type indices for plots instead of drop indices,
#  read drop indices and experiment name from xls file"
#
# The "Niv Experiment" placeholder that used to fill column B (Experiment
# type) for every row is replaced with the real experiment-type labels
# read off the sheet: rows 2-4 belong to the "5uM Las17" experiment and
# rows 5-6 belong to the "5uM Las27" experiment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "5uM Las17"
$ws.Range("B3").Value = "5uM Las17"
$ws.Range("B4").Value = "5uM Las17"
$ws.Range("B5").Value = "5uM Las27"
$ws.Range("B6").Value = "5uM Las27"

# Move the active selection from D8 to B8, matching where the author was
# working (the new Experiment-type column).
$ws.Range("B8").Select() | Out-Null
